# OrderingSuite.xlsx update — "Made changes to update test cases"
$wb = $excel.ActiveWorkbook

# --- Sheet "TestCases" ---
$ws1 = $wb.Worksheets.Item("TestCases")
$ws1.Range("A4").Value = "AddPrivateLineTest"
$ws1.Range("A4").Select()

# --- Sheet "Data" ---
$ws2 = $wb.Worksheets.Item("Data")

# Drop the "Expected_Result" column header / values from the LoginTest block
$ws2.Range("F2").ClearContents()
$ws2.Range("F3").ClearContents()
$ws2.Range("F4").ClearContents()
$ws2.Range("F5").ClearContents()

# LoginTest data rows: flip Runmode to Y and swap in the new username/password pairs
$ws2.Range("A3").Value = "Y"
$ws2.Range("D3").Value = "girish"
$ws2.Range("E3").Value = "'1234"

$ws2.Range("A4").Value = "Y"
$ws2.Range("D4").Value = "Pratap"
$ws2.Range("E4").Value = "'6440904"

$ws2.Range("D5").Value = "24769_selenium"

# "AddTest" section becomes "AddPrivateLineTest"
$ws2.Range("A13").Value = "AddPrivateLineTest"
$ws2.Range("D15").Value = "'26/01/2017"
$ws2.Range("E15").Value = "'67369192"

# Build the new Test5 header block in row 21 (copied formatting from the old row-23 header)
$ws2.Range("F23:H23").Copy() | Out-Null
$ws2.Range("F21:H21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws2.Range("F21").Value = "Col4"
$ws2.Range("G21").Value = "Col5"
$ws2.Range("H21").Value = "Col6"

# Row 22 picks up the values that used to live in row 24's F:H columns
$ws2.Range("F22").Value = "C43"
$ws2.Range("G22").Value = "C53"
$ws2.Range("H22").Value = "C63"

# Row 23's F:H now hold what used to be row 25's values, with the header formatting cleared
$ws2.Range("F23:H23").Style = "Normal"
$ws2.Range("F23").Value = "C45"
$ws2.Range("G23").Value = "C54"
$ws2.Range("H23").Value = "C64"

# Row 24's F:H now hold what used to be row 26's values
$ws2.Range("F24").Value = "C46"
$ws2.Range("G24").Value = "C55"
$ws2.Range("H24").Value = "C65"

# Rows 25 and 26 lose their F:H values entirely
$ws2.Range("F25:H25").ClearContents()
$ws2.Range("F26:H26").ClearContents()

$ws2.Range("D16").Select()
